# Atualizado por script em 07-11-2023 08:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 10 and 11 had their match data (columns F:V) swapped.
#    Columns A:E (index/pais/torneio/temporada/data_partida) stay as-is.
# ---------------------------------------------------------------------------

function Get-RowValues($row) {
    $vals = @{}
    for ($c = 6; $c -le 22; $c++) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

$row10 = Get-RowValues 10
$row11 = Get-RowValues 11

for ($c = 6; $c -le 22; $c++) {
    $ws.Cells.Item(10, $c).Value = $row11[$c]
    $ws.Cells.Item(11, $c).Value = $row10[$c]
}

# ---------------------------------------------------------------------------
# 2) A new match (row 90) was appended at the bottom of the sheet.
#    Copy formatting from the last existing row (89) first, then fill values.
# ---------------------------------------------------------------------------

$ws.Range("A89:V89").Copy()
$ws.Range("A90").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = "denmark"
$ws.Cells.Item(90, 3).Value = "1st-division"
$ws.Cells.Item(90, 4).Value = "2023-2024"
$ws.Cells.Item(90, 5).Value = 45236.79166666666
$ws.Cells.Item(90, 6).Value = "Hillerod"
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = "Naestved"
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 1.95
$ws.Cells.Item(90, 11).Value = "30/10/2023 09:42"
$ws.Cells.Item(90, 12).Value = 1.97
$ws.Cells.Item(90, 13).Value = "06/11/2023 18:58"
$ws.Cells.Item(90, 14).Value = 3.73
$ws.Cells.Item(90, 15).Value = "30/10/2023 09:42"
$ws.Cells.Item(90, 16).Value = 3.7
$ws.Cells.Item(90, 17).Value = "06/11/2023 18:58"
$ws.Cells.Item(90, 18).Value = 3.61
$ws.Cells.Item(90, 19).Value = "30/10/2023 09:42"
$ws.Cells.Item(90, 20).Value = 3.73
$ws.Cells.Item(90, 21).Value = "06/11/2023 18:58"
$ws.Cells.Item(90, 22).Value = "https://www.betexplorer.com/football/denmark/1st-division/hillerod-naestved-if/0fN8moTh/"
